$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H ("Absent") set to 1 for all date rows 3-18, except rows 13 and 17
# (which instead get D and E set to 1, since the student attended those dates).
$absentRows = 3..18 | Where-Object { $_ -ne 13 -and $_ -ne 17 }
foreach ($r in $absentRows) {
    $ws.Cells.Item($r, 8).Value = 1
}

# Rows 13 and 17: student attended -> Total Attendance Count (D) and Real (E) = 1
foreach ($r in @(13, 17)) {
    $ws.Cells.Item($r, 4).Value = 1
    $ws.Cells.Item($r, 5).Value = 1
}
